# This workbook keeps a rolling daily price log. This edit adds three new
# observation rows (new survey date 44476) at the top of the data block
# (rows 563-565), pushing all subsequent rows down by three positions.
# The sheet's used range therefore grows from A1:T669 to A1:T672.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 563; everything from old row 563 onward
# (including the previously-last row, 669) shifts down to 566..672.
$ws.Rows("563:565").Insert()

# --- New row 563 ---
$ws.Range("A563").Value = 6
$ws.Range("B563").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C563").Value = "Metropolitana"
$ws.Range("D563").Value = 44476
$ws.Range("E563").Value = 13
$ws.Range("F563").Value = "Fruta"
$ws.Range("G563").Value = 100108
$ws.Range("H563").Value = "Tropicales y subtropicales"
$ws.Range("I563").Value = 100108006
$ws.Range("J563").Value = "Plátano"
$ws.Range("K563").Value = "Sin especificar"
$ws.Range("L563").Value = "Pintón"
$ws.Range("M563").Value = 950
$ws.Range("N563").Value = 19000
$ws.Range("O563").Value = 20000
$ws.Range("P563").Value = 19421
$ws.Range("Q563").Value = "`$/caja 20 kilos"
$ws.Range("R563").Value = "Ecuador"
$ws.Range("S563").Value = 971
$ws.Range("T563").Value = 20

# --- New row 564 ---
$ws.Range("A564").Value = 6
$ws.Range("B564").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C564").Value = "Metropolitana"
$ws.Range("D564").Value = 44476
$ws.Range("E564").Value = 13
$ws.Range("F564").Value = "Fruta"
$ws.Range("G564").Value = 100108
$ws.Range("H564").Value = "Tropicales y subtropicales"
$ws.Range("I564").Value = 100108006
$ws.Range("J564").Value = "Plátano"
$ws.Range("K564").Value = "Sin especificar"
$ws.Range("L564").Value = "Primera Pintón"
$ws.Range("M564").Value = 2340
$ws.Range("N564").Value = 21000
$ws.Range("O564").Value = 22000
$ws.Range("P564").Value = 21530
$ws.Range("Q564").Value = "`$/caja 20 kilos"
$ws.Range("R564").Value = "Ecuador"
$ws.Range("S564").Value = 1076
$ws.Range("T564").Value = 20

# --- New row 565 ---
$ws.Range("A565").Value = 6
$ws.Range("B565").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C565").Value = "Metropolitana"
$ws.Range("D565").Value = 44476
$ws.Range("E565").Value = 13
$ws.Range("F565").Value = "Fruta"
$ws.Range("G565").Value = 100108
$ws.Range("H565").Value = "Tropicales y subtropicales"
$ws.Range("I565").Value = 100108006
$ws.Range("J565").Value = "Plátano"
$ws.Range("K565").Value = "Sin especificar"
$ws.Range("L565").Value = "Primera Verde"
$ws.Range("M565").Value = 960
$ws.Range("N565").Value = 21000
$ws.Range("O565").Value = 22000
$ws.Range("P565").Value = 21583
$ws.Range("Q565").Value = "`$/caja 20 kilos"
$ws.Range("R565").Value = "Ecuador"
$ws.Range("S565").Value = 1079
$ws.Range("T565").Value = 20
